# Insert a new weekly data row for Mora (Blackberry) at row 39.
# All existing rows from 39 downward shift down by one (39 -> 40, ... 75 -> 76).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 39, pushing rows 39..75 down to 40..76.
$ws.Rows("39").Insert()

# Populate the newly inserted row 39 with the new day's record.
$ws.Range("A39").Value = 6
$ws.Range("B39").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C39").Value = "Metropolitana"
$ws.Range("D39").Value = Get-Date -Year 2022 -Month 1 -Day 26 -Hour 0 -Minute 0 -Second 0
$ws.Range("D39").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E39").Value = 13
$ws.Range("F39").Value = "Fruta"
$ws.Range("G39").Value = 100101
$ws.Range("H39").Value = "Berries"
$ws.Range("I39").Value = 100101008
$ws.Range("J39").Value = "Mora"
$ws.Range("K39").Value = "Sin especificar"
$ws.Range("L39").Value = "Primera"
$ws.Range("M39").Value = 250
$ws.Range("N39").Value = 6000
$ws.Range("O39").Value = 6000
$ws.Range("P39").Value = 6000
$ws.Range("Q39").Value = "$/bandeja 2 kilos"
$ws.Range("R39").Value = "Provincia de Linares"
$ws.Range("S39").Value = 3000
$ws.Range("T39").Value = 2
